# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded for the 37e5cb7e-861c-40ec-816c-c1383e08f148 file across the
# Overview / zh-cn / de-de report sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 37e5cb7e row (row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-12 06:59:31"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the 37e5cb7e row (row 2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-12 06:59:25"
$wsZhCn.Range("K2").Value = "2016-08-12 06:59:52"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the 37e5cb7e row (row 2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-12 06:59:31"
$wsDeDe.Range("K2").Value = "2016-08-12 07:00:07"
